# Updating overlay function to enable users to specify icons for each
# behavioral event: add a new "Overlay" worksheet at the end of the
# workbook with header columns UID / Frame / Behavior / Icon.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the
# end of the tab strip (sheetId 4 / rId4), matching "Overlay" becoming
# the 4th, active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "Overlay"

# Header row for the new Overlay tab.
$ws.Range("A1").Value = "UID"
$ws.Range("B1").Value = "Frame"
$ws.Range("C1").Value = "Behavior"
$ws.Range("D1").Value = "Icon"

# Match the author's final selection/view state on the new tab.
$ws.Range("G41").Select()
